$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = 22
$ws.Range("B23").Value = 2.3190972222222221
$ws.Range("B23").NumberFormat = "[h]:mm:ss"
$ws.Range("C23").Formula = "=SUM(B2:B23)+1.2708333333"
$ws.Range("C23").NumberFormat = "[h]:mm:ss"
$ws.Range("D23").Value = "¿Quién mató a Sara? (Audiovisual, Spanish, Re-watch):39; Harry Potter book 7 (Text-only, English, Familiar):42; [LA CIUDAD PERFECTA - El Futuro de Stephen Hawking - Documental 720p](https://youtu.be/k7nQSaM5VjE) (Audiovisual, English, New):43; Squid Game (Audiovisual, Korean, New):37;"

$ws.Range("C24").Select()
